$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Cells.Item(2, 12).Value = 4611
$ws.Cells.Item(3, 12).Value = 4963
$ws.Cells.Item(4, 12).Value = 1226
$ws.Cells.Item(5, 12).Value = 287
$ws.Cells.Item(6, 12).Value = 4232
$ws.Cells.Item(7, 12).Value = 15319

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Cells.Item(4, 12).Value = 59
$ws.Cells.Item(7, 12).Value = 505
$ws.Cells.Item(8, 12).Value = 1018
$ws.Cells.Item(9, 12).Value = 92
$ws.Cells.Item(11, 12).Value = 246
$ws.Cells.Item(13, 12).Value = 25
$ws.Cells.Item(14, 12).Value = 77
$ws.Cells.Item(16, 12).Value = 33
$ws.Cells.Item(18, 12).Value = 108
$ws.Cells.Item(19, 12).Value = 421
$ws.Cells.Item(20, 12).Value = 389
$ws.Cells.Item(25, 12).Value = 89
$ws.Cells.Item(29, 12).Value = 849
$ws.Cells.Item(31, 12).Value = 149
$ws.Cells.Item(33, 12).Value = 697
$ws.Cells.Item(36, 12).Value = 200
$ws.Cells.Item(37, 12).Value = 564
$ws.Cells.Item(42, 12).Value = 497
$ws.Cells.Item(43, 12).Value = 111
$ws.Cells.Item(44, 12).Value = 109
$ws.Cells.Item(46, 12).Value = 36
$ws.Cells.Item(48, 12).Value = 198
$ws.Cells.Item(49, 12).Value = 80
$ws.Cells.Item(51, 12).Value = 189
$ws.Cells.Item(52, 12).Value = 310
$ws.Cells.Item(53, 12).Value = 177
$ws.Cells.Item(54, 12).Value = 320
$ws.Cells.Item(60, 12).Value = 96
$ws.Cells.Item(63, 12).Value = 43
$ws.Cells.Item(65, 12).Value = 297
$ws.Cells.Item(66, 12).Value = 38
$ws.Cells.Item(67, 12).Value = 528
$ws.Cells.Item(68, 12).Value = 50
$ws.Cells.Item(71, 12).Value = 45
$ws.Cells.Item(76, 12).Value = 239
$ws.Cells.Item(78, 12).Value = 201
$ws.Cells.Item(79, 12).Value = 405
$ws.Cells.Item(82, 12).Value = 22
$ws.Cells.Item(83, 12).Value = 334
$ws.Cells.Item(84, 12).Value = 148
$ws.Cells.Item(85, 12).Value = 782
$ws.Cells.Item(87, 12).Value = 44
$ws.Cells.Item(90, 12).Value = 154
$ws.Cells.Item(94, 12).Value = 192
$ws.Cells.Item(95, 12).Value = 207
$ws.Cells.Item(96, 12).Value = 173
$ws.Cells.Item(97, 12).Value = 131
$ws.Cells.Item(101, 12).Value = 15319

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Cells.Item(3, 12).Value = 16
$ws.Cells.Item(7, 12).Value = 77

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Cells.Item(6, 12).Value = 50
$ws.Cells.Item(7, 12).Value = 173

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Cells.Item(2, 12).Value = 170
$ws.Cells.Item(3, 12).Value = 167
$ws.Cells.Item(7, 12).Value = 505

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Cells.Item(2, 12).Value = 91
$ws.Cells.Item(6, 12).Value = 58
$ws.Cells.Item(7, 12).Value = 246

$ws = $wb.Worksheets.Item("South Shore")
$ws.Cells.Item(3, 12).Value = 315
$ws.Cells.Item(7, 12).Value = 782

$ws = $wb.Worksheets.Item("Little Village")
$ws.Cells.Item(2, 12).Value = 104
$ws.Cells.Item(6, 12).Value = 83
$ws.Cells.Item(7, 12).Value = 310

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Cells.Item(3, 12).Value = 44
$ws.Cells.Item(6, 12).Value = 60
$ws.Cells.Item(7, 12).Value = 177

$ws = $wb.Worksheets.Item("Austin")
$ws.Cells.Item(3, 12).Value = 344
$ws.Cells.Item(4, 12).Value = 78
$ws.Cells.Item(6, 12).Value = 269
$ws.Cells.Item(7, 12).Value = 1018

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Cells.Item(2, 12).Value = 105
$ws.Cells.Item(3, 12).Value = 130
$ws.Cells.Item(6, 12).Value = 78
$ws.Cells.Item(7, 12).Value = 334

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Cells.Item(2, 12).Value = 190
$ws.Cells.Item(3, 12).Value = 236
$ws.Cells.Item(5, 12).Value = 15
$ws.Cells.Item(7, 12).Value = 697

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Cells.Item(2, 12).Value = 79
$ws.Cells.Item(7, 12).Value = 207

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Cells.Item(2, 12).Value = 168
$ws.Cells.Item(3, 12).Value = 196
$ws.Cells.Item(7, 12).Value = 564

$ws = $wb.Worksheets.Item("New City")
$ws.Cells.Item(3, 12).Value = 94
$ws.Cells.Item(7, 12).Value = 297

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Cells.Item(3, 12).Value = 37
$ws.Cells.Item(6, 12).Value = 43
$ws.Cells.Item(7, 12).Value = 149

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Cells.Item(2, 12).Value = 153
$ws.Cells.Item(3, 12).Value = 205
$ws.Cells.Item(4, 12).Value = 38
$ws.Cells.Item(7, 12).Value = 528

$ws = $wb.Worksheets.Item("South Deering")
$ws.Cells.Item(6, 12).Value = 41
$ws.Cells.Item(7, 12).Value = 148

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Cells.Item(2, 12).Value = 28
$ws.Cells.Item(7, 12).Value = 80

$ws = $wb.Worksheets.Item("Loop")
$ws.Cells.Item(6, 12).Value = 156
$ws.Cells.Item(7, 12).Value = 320

$ws = $wb.Worksheets.Item("Englewood")
$ws.Cells.Item(2, 12).Value = 256
$ws.Cells.Item(3, 12).Value = 320
$ws.Cells.Item(7, 12).Value = 849

$ws = $wb.Worksheets.Item("Lake View")
$ws.Cells.Item(4, 12).Value = 39
$ws.Cells.Item(7, 12).Value = 198

$ws = $wb.Worksheets.Item("Chatham")
$ws.Cells.Item(2, 12).Value = 149
$ws.Cells.Item(3, 12).Value = 130
$ws.Cells.Item(7, 12).Value = 421

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Cells.Item(3, 12).Value = 30
$ws.Cells.Item(7, 12).Value = 109

$ws = $wb.Worksheets.Item("River North")
$ws.Cells.Item(3, 12).Value = 43
$ws.Cells.Item(6, 12).Value = 112
$ws.Cells.Item(7, 12).Value = 239

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Cells.Item(2, 12).Value = 142
$ws.Cells.Item(3, 12).Value = 168
$ws.Cells.Item(6, 12).Value = 139
$ws.Cells.Item(7, 12).Value = 497

$ws = $wb.Worksheets.Item("Boystown")
$ws.Cells.Item(3, 12).Value = 8
$ws.Cells.Item(6, 12).Value = 25

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Cells.Item(3, 12).Value = 65
$ws.Cells.Item(7, 12).Value = 201

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Cells.Item(6, 12).Value = 12
$ws.Cells.Item(7, 12).Value = 36

$ws = $wb.Worksheets.Item("Roseland")
$ws.Cells.Item(3, 12).Value = 147
$ws.Cells.Item(5, 12).Value = 12
$ws.Cells.Item(7, 12).Value = 405

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Cells.Item(3, 12).Value = 128
$ws.Cells.Item(7, 12).Value = 389

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Cells.Item(3, 12).Value = 39
$ws.Cells.Item(7, 12).Value = 108

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Cells.Item(3, 12).Value = 59
$ws.Cells.Item(7, 12).Value = 200

$ws = $wb.Worksheets.Item("West Loop")
$ws.Cells.Item(2, 12).Value = 44
$ws.Cells.Item(7, 12).Value = 192

$ws = $wb.Worksheets.Item("East Side")
$ws.Cells.Item(2, 12).Value = 33
$ws.Cells.Item(6, 12).Value = 12
$ws.Cells.Item(7, 12).Value = 89

$ws = $wb.Worksheets.Item("North Center")
$ws.Cells.Item(2, 12).Value = 10
$ws.Cells.Item(3, 12).Value = 8
$ws.Cells.Item(7, 12).Value = 38

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Cells.Item(6, 12).Value = 25
$ws.Cells.Item(7, 12).Value = 92

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Cells.Item(3, 12).Value = 37
$ws.Cells.Item(6, 12).Value = 29

$ws = $wb.Worksheets.Item("West Town")
$ws.Cells.Item(6, 12).Value = 66
$ws.Cells.Item(7, 12).Value = 131

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Cells.Item(6, 12).Value = 43
$ws.Cells.Item(7, 12).Value = 154

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Cells.Item(3, 12).Value = 63
$ws.Cells.Item(7, 12).Value = 189

$ws = $wb.Worksheets.Item("North Park")
$ws.Cells.Item(3, 12).Value = 16
$ws.Cells.Item(7, 12).Value = 50

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Cells.Item(6, 12).Value = 25
$ws.Cells.Item(7, 12).Value = 96

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Cells.Item(3, 12).Value = 34
$ws.Cells.Item(7, 12).Value = 111

$ws = $wb.Worksheets.Item("Oakland")
$ws.Cells.Item(2, 12).Value = 19
$ws.Cells.Item(3, 12).Value = 16
$ws.Cells.Item(7, 12).Value = 45

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Cells.Item(6, 12).Value = 7
$ws.Cells.Item(7, 12).Value = 22

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Cells.Item(6, 12).Value = 19
$ws.Cells.Item(7, 12).Value = 59

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Cells.Item(3, 12).Value = 10
$ws.Cells.Item(4, 12).Value = 6
$ws.Cells.Item(7, 12).Value = 44

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Cells.Item(2, 12).Value = 5
$ws.Cells.Item(7, 12).Value = 33
